$d = $word.ActiveDocument

# --- Change 1 -------------------------------------------------------
# Merge the three runs "    hotel_ID SERIAL PRIMARY KEY NOT NULL," + " "
# + "--auto generating ID" into a single run with the combined text.
# Find & Replace (even with an unchanged net string) makes the engine
# rebuild the matched range as one run, collapsing the extra run splits.
$rng1 = $d.Content
$rng1.Find.Execute(
    "    hotel_ID SERIAL PRIMARY KEY NOT NULL, --auto generating ID",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "    hotel_ID SERIAL PRIMARY KEY NOT NULL, --auto generating ID", 2)

# --- Change 2 -------------------------------------------------------
# Turn the empty paragraph right after "MySQL" into a new bulleted
# list item "CSS" that continues the same list (numId 4, ilvl 0) as
# the Python/Javascript/MySQL items above it.
$rng2 = $d.Content
$rng2.Find.Execute("MySQL")
$mysqlPara = $rng2.Paragraphs(1)
$cssPara = $mysqlPara.Next()
$cssPara.Range.Text = "CSS"
$cssPara.Style = $mysqlPara.Style
$cssPara.Range.ListFormat.ApplyListTemplateWithLevel($mysqlPara.Range.ListFormat.ListTemplate, $true)

# --- Change 3 -------------------------------------------------------
# Remove the yellow highlight from the paragraph mark of the
# "create view rooms_avail" bullet, and delete the trailing
# "[do I expand to include entire statement?]" run entirely (its own
# yellow-highlighted run disappears along with it).
$rng3 = $d.Content
$rng3.Find.Execute("create view rooms_avail")
$viewPara = $rng3.Paragraphs(1)
$capacityPara = $viewPara.Next()
$viewPara.Style = $capacityPara.Style
$viewPara.Range.ListFormat.ApplyListTemplateWithLevel($capacityPara.Range.ListFormat.ListTemplate, $true)
$viewPara.Range.ListFormat.ListLevelNumber = 2

$rng4 = $d.Content
$rng4.Find.Execute("[do I expand to include entire statement?]")
$rng4.Text = ""
